$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 111618144
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "2"
$ws.Range("I4").Style = "Normal"
$ws.Range("K4").Value = ""
$ws.Range("Q4").Value = 580620.6996611424
$ws.Range("R4").Value = 6415142.541277731
$ws.Range("AC4").Value = ""

# Row 5
$ws.Range("A5").Value = 111618046
$ws.Range("B5").Value = 93388
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 2180
$ws.Range("F5").Value = "Blåmossa"
$ws.Range("G5").Value = "Leucobryum glaucum"
$ws.Range("H5").Value = "(Hedw.) Ångstr."
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("Q5").Value = 580591.6383206119
$ws.Range("R5").Value = 6415156.322361182
$ws.Range("AC5").Value = ""

# Row 6
$ws.Range("A6").Value = 111618070
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "15"
$ws.Range("I6").Style = "Normal"
$ws.Range("K6").Value = "blomning"
$ws.Range("Q6").Value = 580592.470229132
$ws.Range("R6").Value = 6415141.442167919
$ws.Range("AC6").Value = "1 blomma"

# Row 7
$ws.Range("A7").Value = 111618078
$ws.Range("B7").Value = 96348
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = "Knärot"
$ws.Range("G7").Value = "Goodyera repens"
$ws.Range("H7").Value = "(L.) R. Br."
$ws.Range("J7").Value = "plantor/tuvor"
$ws.Range("Q7").Value = 580612.1009209087
$ws.Range("R7").Value = 6415119.491031807

# Row 8
$ws.Range("A8").Value = 111618089
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "30"
$ws.Range("I8").Style = "Normal"
$ws.Range("P8").Value = "A 32649, Heda, Sm"
$ws.Range("Q8").Value = 580617.6201989455
$ws.Range("R8").Value = 6415136.627037819

# Row 9
$ws.Range("A9").Value = 111618109
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "10"
$ws.Range("I9").Style = "Normal"
$ws.Range("K9").Value = "blomning"
$ws.Range("P9").Value = "A 32649, Sm"
$ws.Range("Q9").Value = 580619.1666838422
$ws.Range("R9").Value = 6415112.716507593
$ws.Range("AC9").Value = "1 blomma"
